$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new pretty-printed JSON-style text that replaces the old A2 value
$newText = 'questions = [
    {
        "title": "The code below contains an error. Which option describes the correct way to fix it? const request = $.get(\"/test\")\nrequest.success((data)=&gt{\n    data = JSON.parse(data)\n    console.log(data)\n})",
        "ques_type": 15,
        "options": [
            "A success function must be passed as a parameter to $.get().",
            "The line data = JSON.parse(data) should be removed.",
            "The function request.success should be replaced with request.then.",
            "The arrow function should be changed to anonymous function."
        ],
        "score": [
            "The function request.success should be replaced with request.then."
        ]
    },
    {
        "title": "In the code below, which property contains the response status code? const request = $.ajax(\"/test\", {\n           type: \"GET\",\n           success: function(data, textStatus, xhr){\n               //Display status code\n           }\n})",
        "ques_type": 2,
        "options": [
            "textStatus",
            "xhr.status",
            "data.statusCode",
            "xhr.code"
        ],
        "score": "xhr.status"
    },
    {
        "title": "True or false: The response from the fetch request will be displayed in the console. let response\nfetch(\"/test\")\n.then(response =&gt response.json())\n.then(data =&gt response = data)\n \nconsole.log(response)",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    },
    {
        "title": "The request below returns the following data: {\"data\": [\n   1,2,3,4\n]}\n You want to insert the four numbers in the response data above as div elements in the container. Which snippet of code accomplishes this? &ltdiv class=\"container\"&gt&lt/div&gt\n   &ltscript&gt\n       let response\n       fetch(\"/test\")\n       .then(response =&gt response.json())\n       .then(data =&gt {\n             //your code goes here\n         })\n   &lt/script&gt",
        "ques_type": 2,
        "options": [
            "const container = document.querySelector(\".container\")\ndata.data.forEach((elm)=&gt{\n   container.append(\"&ltdiv&gt\"+elm+\"&lt/div&gt\")\n})\n",
            "const container = document.querySelector(\".container\")\ndata.data.forEach((elm)=&gt{\n     container.innerHTML += \"&ltdiv&gt\"+elm+\"&lt/div&gt\"\n })\n",
            "data = JSON.parse(data)\nconst container = document.querySelector(\".container\")\ndata.forEach((elm)=&gt{\n   const div = document.createElement(\"div\")\n   div.innerText = elm\n   container.append(div)\n})\n",
            "const container = document.querySelector(\".container\")\ndata.map((elm)=&gt{\n   const div = document.createElement(\"div\")\n   div.innerText = elm\n   container.append(div)\n})\n"
        ],
        "score": "const container = document.querySelector(\".container\")\ndata.data.forEach((elm)=&gt{\n     container.innerHTML += \"&ltdiv&gt\"+elm+\"&lt/div&gt\"\n })"
    }
]'

# Remove the old row 2 (which held the shared-string text) - content moves into A1
$ws.Rows(2).Delete()

# A1 previously held a bordered/bold/centered "0" placeholder value & style;
# strip that formatting so A1 becomes a plain, unstyled cell
$ws.Range("A1").ClearFormats()

# Put the reformatted questions text into A1
$ws.Range("A1").Value = $newText

# Avoid a stale "custom height" row dimension being written out for the
# now much taller multi-line text
$ws.Rows(1).EntireRow.AutoFit()
